$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,20
$arr[0,0] = "ECs"
$arr[0,1] = "Angpt1"
$arr[0,2] = "Itga5"
$arr[0,3] = "ECs"
$arr[0,4] = 1
$arr[0,5] = 0.5
$arr[0,6] = 0.1657605
$arr[0,7] = 0.331521
$arr[0,8] = 0.01219521989678654
$arr[0,9] = 0.008352702481122073
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 42.960745
$arr[0,13] = 85.92149000000001
$arr[0,14] = 0.1929389995390091
$arr[0,15] = 0.1411745935894314
$arr[0,16] = 7.1211945715725
$arr[0,17] = 28.48477828629
$arr[0,18] = 0.002352933526044214
$arr[0,19] = 0.001179189378145844
$arr[1,0] = "ECs"
$arr[1,1] = "Angpt1"
$arr[1,2] = "Itga5"
$arr[1,3] = "FAPs"
$arr[1,4] = 1
$arr[1,5] = 0.5
$arr[1,6] = 0.1657605
$arr[1,7] = 0.331521
$arr[1,8] = 0.01219521989678654
$arr[1,9] = 0.008352702481122073
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 69.97731266666666
$arr[1,13] = 209.931938
$arr[1,14] = 0.3142718473884727
$arr[1,15] = 0.3449318212311228
$arr[1,16] = 11.599474336283
$arr[1,17] = 69.596846017698
$arr[1,18] = 0.003832614286271767
$arr[1,19] = 0.002881112879015154
$arr[2,0] = "ECs"
$arr[2,1] = "Angpt1"
$arr[2,2] = "Itga5"
$arr[2,3] = "Inflammatory-Mac"
$arr[2,4] = 1
$arr[2,5] = 0.5
$arr[2,6] = 0.1657605
$arr[2,7] = 0.331521
$arr[2,8] = 0.01219521989678654
$arr[2,9] = 0.008352702481122073
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 38.33270933333333
$arr[2,13] = 114.998128
$arr[2,14] = 0.1721542442616619
$arr[2,15] = 0.1889494000155887
$arr[2,16] = 6.354049065448001
$arr[2,17] = 38.12429439268801
$arr[2,18] = 0.00209945886493607
$arr[2,19] = 0.001578238122316735
$arr[3,0] = "ECs"
$arr[3,1] = "Angpt1"
$arr[3,2] = "Itga5"
$arr[3,3] = "MuSCs"
$arr[3,4] = 1
$arr[3,5] = 0.5
$arr[3,6] = 0.1657605
$arr[3,7] = 0.331521
$arr[3,8] = 0.01219521989678654
$arr[3,9] = 0.008352702481122073
$arr[3,10] = 2
$arr[3,11] = 1
$arr[3,12] = 16.4153395
$arr[3,13] = 32.830679
$arr[3,14] = 0.07372216613615937
$arr[3,15] = 0.05394293982902391
$arr[3,16] = 2.72101488318975
$arr[3,17] = 10.884059532759
$arr[3,18] = 0.0008990580272978938
$arr[3,19] = 0.0004505693273489067
$arr[4,0] = "ECs"
$arr[4,1] = "Angpt1"
$arr[4,2] = "Itga5"
$arr[4,3] = "Neutrophils"
$arr[4,4] = 1
$arr[4,5] = 0.5
$arr[4,6] = 0.1657605
$arr[4,7] = 0.331521
$arr[4,8] = 0.01219521989678654
$arr[4,9] = 0.008352702481122073
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 26.08950466666667
$arr[4,13] = 78.26851400000001
$arr[4,14] = 0.1171693584190632
$arr[4,15] = 0.1286002565225384
$arr[4,16] = 4.324609338299001
$arr[4,17] = 25.947656029794
$arr[4,18] = 0.001428906091085873
$arr[4,19] = 0.001074159681728742
$arr[5,0] = "ECs"
$arr[5,1] = "Angpt1"
$arr[5,2] = "Itga5"
$arr[5,3] = "Resolving-Mac"
$arr[5,4] = 1
$arr[5,5] = 0.5
$arr[5,6] = 0.1657605
$arr[5,7] = 0.331521
$arr[5,8] = 0.01219521989678654
$arr[5,9] = 0.008352702481122073
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 28.889299
$arr[5,13] = 86.66789700000001
$arr[5,14] = 0.1297433842556338
$arr[5,15] = 0.1424009888122948
$arr[5,16] = 4.788704646889501
$arr[5,17] = 28.732227881337
$arr[5,18] = 0.001582249101150727
$arr[5,19] = 0.001189433092566691
$arr[6,0] = "FAPs"
$arr[6,1] = "Angpt1"
$arr[6,2] = "Itga5"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 12.45740533333333
$arr[6,7] = 37.372216
$arr[6,8] = 0.91650783741242
$arr[6,9] = 0.9415964638989084
$arr[6,10] = 2
$arr[6,11] = 1
$arr[6,12] = 42.960745
$arr[6,13] = 85.92149000000001
$arr[6,14] = 0.1929389995390091
$arr[6,15] = 0.1411745935894314
$arr[6,16] = 535.1794138869734
$arr[6,17] = 3211.07648332184
$arr[6,18] = 0.1768301052200132
$arr[6,19] = 0.1329294981161741
$arr[7,0] = "FAPs"
$arr[7,1] = "Angpt1"
$arr[7,2] = "Itga5"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 12.45740533333333
$arr[7,7] = 37.372216
$arr[7,8] = 0.91650783741242
$arr[7,9] = 0.9415964638989084
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 69.97731266666666
$arr[7,13] = 209.931938
$arr[7,14] = 0.3142718473884727
$arr[7,15] = 0.3449318212311228
$arr[7,16] = 871.7357480260675
$arr[7,17] = 7845.621732234608
$arr[7,18] = 0.2880326112096153
$arr[7,19] = 0.3247865831574356
$arr[8,0] = "FAPs"
$arr[8,1] = "Angpt1"
$arr[8,2] = "Itga5"
$arr[8,3] = "Inflammatory-Mac"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 12.45740533333333
$arr[8,7] = 37.372216
$arr[8,8] = 0.91650783741242
$arr[8,9] = 0.9415964638989084
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 38.33270933333333
$arr[8,13] = 114.998128
$arr[8,14] = 0.1721542442616619
$arr[8,15] = 0.1889494000155887
$arr[8,16] = 477.5260976901831
$arr[8,17] = 4297.734879211648
$arr[8,18] = 0.1577807141096253
$arr[8,19] = 0.1779140869104987
$arr[9,0] = "FAPs"
$arr[9,1] = "Angpt1"
$arr[9,2] = "Itga5"
$arr[9,3] = "MuSCs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 12.45740533333333
$arr[9,7] = 37.372216
$arr[9,8] = 0.91650783741242
$arr[9,9] = 0.9415964638989084
$arr[9,10] = 2
$arr[9,11] = 1
$arr[9,12] = 16.4153395
$arr[9,13] = 32.830679
$arr[9,14] = 0.07372216613615937
$arr[9,15] = 0.05394293982902391
$arr[9,16] = 204.4925378357774
$arr[9,17] = 1226.955227014664
$arr[9,18] = 0.06756694305481056
$arr[9,19] = 0.0507924813953205
$arr[10,0] = "FAPs"
$arr[10,1] = "Angpt1"
$arr[10,2] = "Itga5"
$arr[10,3] = "Neutrophils"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 12.45740533333333
$arr[10,7] = 37.372216
$arr[10,8] = 0.91650783741242
$arr[10,9] = 0.9415964638989084
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 26.08950466666667
$arr[10,13] = 78.26851400000001
$arr[10,14] = 0.1171693584190632
$arr[10,15] = 0.1286002565225384
$arr[10,16] = 325.0075345785583
$arr[10,17] = 2925.067811207025
$arr[10,18] = 0.1073866352956563
$arr[10,19] = 0.1210895467981147
$arr[11,0] = "FAPs"
$arr[11,1] = "Angpt1"
$arr[11,2] = "Itga5"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 12.45740533333333
$arr[11,7] = 37.372216
$arr[11,8] = 0.91650783741242
$arr[11,9] = 0.9415964638989084
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 28.889299
$arr[11,13] = 86.66789700000001
$arr[11,14] = 0.1297433842556338
$arr[11,15] = 0.1424009888122948
$arr[11,16] = 359.8857074388614
$arr[11,17] = 3238.971366949752
$arr[11,18] = 0.1189108285226995
$arr[11,19] = 0.1340842675213648
$arr[12,0] = "MuSCs"
$arr[12,1] = "Angpt1"
$arr[12,2] = "Itga5"
$arr[12,3] = "ECs"
$arr[12,4] = 2
$arr[12,5] = 1
$arr[12,6] = 0.9207270000000001
$arr[12,7] = 1.841454
$arr[12,8] = 0.06773910690368684
$arr[12,9] = 0.04639560508888476
$arr[12,10] = 2
$arr[12,11] = 1
$arr[12,12] = 42.960745
$arr[12,13] = 85.92149000000001
$arr[12,14] = 0.1929389995390091
$arr[12,15] = 0.1411745935894314
$arr[12,16] = 39.555117861615
$arr[12,17] = 158.22047144646
$arr[12,18] = 0.01306951551566332
$arr[12,19] = 0.00654988069275906
$arr[13,0] = "MuSCs"
$arr[13,1] = "Angpt1"
$arr[13,2] = "Itga5"
$arr[13,3] = "FAPs"
$arr[13,4] = 2
$arr[13,5] = 1
$arr[13,6] = 0.9207270000000001
$arr[13,7] = 1.841454
$arr[13,8] = 0.06773910690368684
$arr[13,9] = 0.04639560508888476
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 69.97731266666666
$arr[13,13] = 209.931938
$arr[13,14] = 0.3142718473884727
$arr[13,15] = 0.3449318212311228
$arr[13,16] = 64.430001159642
$arr[13,17] = 386.580006957852
$arr[13,18] = 0.02128849426706691
$arr[13,19] = 0.01600332056042897
$arr[14,0] = "MuSCs"
$arr[14,1] = "Angpt1"
$arr[14,2] = "Itga5"
$arr[14,3] = "Inflammatory-Mac"
$arr[14,4] = 2
$arr[14,5] = 1
$arr[14,6] = 0.9207270000000001
$arr[14,7] = 1.841454
$arr[14,8] = 0.06773910690368684
$arr[14,9] = 0.04639560508888476
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 38.33270933333333
$arr[14,13] = 114.998128
$arr[14,14] = 0.1721542442616619
$arr[14,15] = 0.1889494000155887
$arr[14,16] = 35.293960466352
$arr[14,17] = 211.763762798112
$arr[14,18] = 0.01166157475596413
$arr[14,19] = 0.008766421744904971
$arr[15,0] = "MuSCs"
$arr[15,1] = "Angpt1"
$arr[15,2] = "Itga5"
$arr[15,3] = "MuSCs"
$arr[15,4] = 2
$arr[15,5] = 1
$arr[15,6] = 0.9207270000000001
$arr[15,7] = 1.841454
$arr[15,8] = 0.06773910690368684
$arr[15,9] = 0.04639560508888476
$arr[15,10] = 2
$arr[15,11] = 1
$arr[15,12] = 16.4153395
$arr[15,13] = 32.830679
$arr[15,14] = 0.07372216613615937
$arr[15,15] = 0.05394293982902391
$arr[15,16] = 15.1140462918165
$arr[15,17] = 60.45618516726601
$arr[15,18] = 0.004993873693068661
$arr[15,19] = 0.002502715333640866
$arr[16,0] = "MuSCs"
$arr[16,1] = "Angpt1"
$arr[16,2] = "Itga5"
$arr[16,3] = "Neutrophils"
$arr[16,4] = 2
$arr[16,5] = 1
$arr[16,6] = 0.9207270000000001
$arr[16,7] = 1.841454
$arr[16,8] = 0.06773910690368684
$arr[16,9] = 0.04639560508888476
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = 26.08950466666667
$arr[16,13] = 78.26851400000001
$arr[16,14] = 0.1171693584190632
$arr[16,15] = 0.1286002565225384
$arr[16,16] = 24.021311363226
$arr[16,17] = 144.127868179356
$arr[16,18] = 0.00793694769578532
$arr[16,19] = 0.005966486715948969
$arr[17,0] = "MuSCs"
$arr[17,1] = "Angpt1"
$arr[17,2] = "Itga5"
$arr[17,3] = "Resolving-Mac"
$arr[17,4] = 2
$arr[17,5] = 1
$arr[17,6] = 0.9207270000000001
$arr[17,7] = 1.841454
$arr[17,8] = 0.06773910690368684
$arr[17,9] = 0.04639560508888476
$arr[17,10] = 3
$arr[17,11] = 1
$arr[17,12] = 28.889299
$arr[17,13] = 86.66789700000001
$arr[17,14] = 0.1297433842556338
$arr[17,15] = 0.1424009888122948
$arr[17,16] = 26.59915760037301
$arr[17,17] = 159.594945602238
$arr[17,18] = 0.008788700976138495
$arr[17,19] = 0.006606780041201926
$arr[18,0] = "Neutrophils"
$arr[18,1] = "Angpt1"
$arr[18,2] = "Itga5"
$arr[18,3] = "ECs"
$arr[18,4] = 1
$arr[18,5] = 0.3333333333333333
$arr[18,6] = 0.04835900000000001
$arr[18,7] = 0.145077
$arr[18,8] = 0.003557835787106702
$arr[18,9] = 0.003655228531084749
$arr[18,10] = 2
$arr[18,11] = 1
$arr[18,12] = 42.960745
$arr[18,13] = 85.92149000000001
$arr[18,14] = 0.1929389995390091
$arr[18,15] = 0.1411745935894314
$arr[18,16] = 2.077538667455
$arr[18,17] = 12.46523200473
$arr[18,18] = 0.0006864452772884501
$arr[18,19] = 0.0005160254023523836
$arr[19,0] = "Neutrophils"
$arr[19,1] = "Angpt1"
$arr[19,2] = "Itga5"
$arr[19,3] = "FAPs"
$arr[19,4] = 1
$arr[19,5] = 0.3333333333333333
$arr[19,6] = 0.04835900000000001
$arr[19,7] = 0.145077
$arr[19,8] = 0.003557835787106702
$arr[19,9] = 0.003655228531084749
$arr[19,10] = 3
$arr[19,11] = 1
$arr[19,12] = 69.97731266666666
$arr[19,13] = 209.931938
$arr[19,14] = 0.3142718473884727
$arr[19,15] = 0.3449318212311228
$arr[19,16] = 3.384032863247334
$arr[19,17] = 30.456295769226
$arr[19,18] = 0.001118127625518844
$arr[19,19] = 0.001260804634243024
$arr[20,0] = "Neutrophils"
$arr[20,1] = "Angpt1"
$arr[20,2] = "Itga5"
$arr[20,3] = "Inflammatory-Mac"
$arr[20,4] = 1
$arr[20,5] = 0.3333333333333333
$arr[20,6] = 0.04835900000000001
$arr[20,7] = 0.145077
$arr[20,8] = 0.003557835787106702
$arr[20,9] = 0.003655228531084749
$arr[20,10] = 3
$arr[20,11] = 1
$arr[20,12] = 38.33270933333333
$arr[20,13] = 114.998128
$arr[20,14] = 0.1721542442616619
$arr[20,15] = 0.1889494000155887
$arr[20,16] = 1.853731490650667
$arr[20,17] = 16.683583415856
$arr[20,18] = 0.0006124965311364493
$arr[20,19] = 0.000690653237868325
$arr[21,0] = "Neutrophils"
$arr[21,1] = "Angpt1"
$arr[21,2] = "Itga5"
$arr[21,3] = "MuSCs"
$arr[21,4] = 1
$arr[21,5] = 0.3333333333333333
$arr[21,6] = 0.04835900000000001
$arr[21,7] = 0.145077
$arr[21,8] = 0.003557835787106702
$arr[21,9] = 0.003655228531084749
$arr[21,10] = 2
$arr[21,11] = 1
$arr[21,12] = 16.4153395
$arr[21,13] = 32.830679
$arr[21,14] = 0.07372216613615937
$arr[21,15] = 0.05394293982902391
$arr[21,16] = 0.7938294028805002
$arr[21,17] = 4.762976417283001
$arr[21,18] = 0.0002622913609822536
$arr[21,19] = 0.000197173772713636
$arr[22,0] = "Neutrophils"
$arr[22,1] = "Angpt1"
$arr[22,2] = "Itga5"
$arr[22,3] = "Neutrophils"
$arr[22,4] = 1
$arr[22,5] = 0.3333333333333333
$arr[22,6] = 0.04835900000000001
$arr[22,7] = 0.145077
$arr[22,8] = 0.003557835787106702
$arr[22,9] = 0.003655228531084749
$arr[22,10] = 3
$arr[22,11] = 1
$arr[22,12] = 26.08950466666667
$arr[22,13] = 78.26851400000001
$arr[22,14] = 0.1171693584190632
$arr[22,15] = 0.1286002565225384
$arr[22,16] = 1.261662356175334
$arr[22,17] = 11.354961205578
$arr[22,18] = 0.0004168693365356749
$arr[22,19] = 0.000470063326746
$arr[23,0] = "Neutrophils"
$arr[23,1] = "Angpt1"
$arr[23,2] = "Itga5"
$arr[23,3] = "Resolving-Mac"
$arr[23,4] = 1
$arr[23,5] = 0.3333333333333333
$arr[23,6] = 0.04835900000000001
$arr[23,7] = 0.145077
$arr[23,8] = 0.003557835787106702
$arr[23,9] = 0.003655228531084749
$arr[23,10] = 3
$arr[23,11] = 1
$arr[23,12] = 28.889299
$arr[23,13] = 86.66789700000001
$arr[23,14] = 0.1297433842556338
$arr[23,15] = 0.1424009888122948
$arr[23,16] = 1.397057610341
$arr[23,17] = 12.573518493069
$arr[23,18] = 0.00046160565564503
$arr[23,19] = 0.00052050815716138
$ws.Range("A2:T25").Value = $arr
